# Update the cryptos table with the latest scraped price/volume figures.
# (Two rows - RenderToken/PancakeSwap and Cosmos/SuiNetwork - also swapped rank order.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a value to a cell while preserving it as plain text even when it
# looks like a number (prices are stored as text in this sheet, e.g. "0.990", "1.90").
function Set-TextCell($range, [string]$text) {
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $ws.Range($range).Value = "'" + $text
    } else {
        $ws.Range($range).Value = $text
    }
}

Set-TextCell 'D2' '59.982.22'
Set-TextCell 'E2' '  +2.67%  '
Set-TextCell 'D3' '3.200.13'
Set-TextCell 'E3' '  +1.40%  '
Set-TextCell 'E4' '  +0.03%  '
Set-TextCell 'D5' '536.62'
Set-TextCell 'E5' '  -0.13%  '
Set-TextCell 'D6' '145.25'
Set-TextCell 'E6' '  +3.96%  '
Set-TextCell 'E7' '  +0.07%  '
Set-TextCell 'E8' '  +3.61%  '
Set-TextCell 'E9' '  +0.65%  '
Set-TextCell 'E10' '  +3.40%  '
Set-TextCell 'D11' '0.434'
Set-TextCell 'E11' '  +3.12%  '
Set-TextCell 'D12' '3.753.62'
Set-TextCell 'E12' '  +1.50%  '
Set-TextCell 'E13' '  -1.10%  '
Set-TextCell 'D14' '26.04'
Set-TextCell 'E14' '  +0.85%  '
Set-TextCell 'E15' '  +2.08%  '
Set-TextCell 'D16' '60.072.60'
Set-TextCell 'E16' '  +2.71%  '
Set-TextCell 'D17' '3.200.22'
Set-TextCell 'E17' '  +1.01%  '
Set-TextCell 'D18' '6.24'
Set-TextCell 'E18' '  +0.37%  '
Set-TextCell 'D19' '13.18'
Set-TextCell 'E19' '  +1.06%  '
Set-TextCell 'D20' '8.30'
Set-TextCell 'E20' '  +0.34%  '
Set-TextCell 'D21' '379.45'
Set-TextCell 'E21' '  +0.59%  '
Set-TextCell 'E22' '  +0.01%  '
Set-TextCell 'E23' '  +1.83%  '
Set-TextCell 'D24' '70.12'
Set-TextCell 'E24' '  -0.09%  '
Set-TextCell 'D25' '8.81'
Set-TextCell 'E25' '  +8.75%  '
Set-TextCell 'E26' '  +1.27%  '
Set-TextCell 'E27' '  -0.24%  '
Set-TextCell 'D28' '0.0₃0896'
Set-TextCell 'E28' '  +3.02%  '
Set-TextCell 'B29' 'PancakeSwap'
Set-TextCell 'C29' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 'D29' '1.90'
Set-TextCell 'E29' '  +0.58%  '
Set-TextCell 'B30' 'RenderToken'
Set-TextCell 'C30' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D30' '6.19'
Set-TextCell 'E30' '  +0.49%  '
Set-TextCell 'D31' '22.35'
Set-TextCell 'E31' '  +1.96%  '
Set-TextCell 'D32' '5.43'
Set-TextCell 'E32' '  +5.04%  '
Set-TextCell 'D33' '1.21'
Set-TextCell 'E33' '  +2.49%  '
Set-TextCell 'D34' '6.68'
Set-TextCell 'E34' '  +6.99%  '
Set-TextCell 'D35' '157.14'
Set-TextCell 'E35' '  -2.21%  '
Set-TextCell 'D36' '1.35'
Set-TextCell 'E36' '  -1.28%  '
Set-TextCell 'D37' '2.798.39'
Set-TextCell 'E37' '  +5.33%  '
Set-TextCell 'E38' '  +0.70%  '
Set-TextCell 'E39' '  +3.47%  '
Set-TextCell 'E40' '  +0.45%  '
Set-TextCell 'E41' '  +0.77%  '
Set-TextCell 'D42' '39.79'
Set-TextCell 'E42' '  +2.87%  '
Set-TextCell 'D43' '0.0293'
Set-TextCell 'E43' '  +4.40%  '
Set-TextCell 'D44' '0.717'
Set-TextCell 'E44' '  +1.31%  '
Set-TextCell 'E45' '  +2.88%  '
Set-TextCell 'D46' '3.244.43'
Set-TextCell 'E46' '  +1.44%  '
Set-TextCell 'D47' '0.990'
Set-TextCell 'E47' '  +1.14%  '
Set-TextCell 'B48' 'SuiNetwork'
Set-TextCell 'C48' 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextCell 'D48' '0.813'
Set-TextCell 'E48' '  +7.28%  '
Set-TextCell 'B49' 'Cosmos'
Set-TextCell 'C49' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 'D49' '6.16'
Set-TextCell 'E49' '  -0.95%  '
Set-TextCell 'D50' '20.61'
Set-TextCell 'E50' '  +1.64%  '
Set-TextCell 'E51' '  +0.01%  '
